$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Bgn","Fgfr3","ECs",3,1,38.718679,116.156037,0.01404461724059496,0.01404461724059496,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,101.4411321083133,912.97018897482,0.008567258002214266,0.008567258002214266),
    @("ECs","Bgn","Fgfr3","FAPs",3,1,38.718679,116.156037,0.01404461724059496,0.01404461724059496,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,22.06377469701833,198.573972273165,0.001863406356015933,0.001863406356015933),
    @("ECs","Bgn","Fgfr3","sCs",3,1,38.718679,116.156037,0.01404461724059496,0.01404461724059496,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,42.791225813257,385.121032319313,0.003613952882364764,0.003613952882364764),
    @("FAPs","Bgn","Fgfr3","ECs",3,1,2511.398112333333,7534.194336999999,0.9109718135392577,0.9109718135392579,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,6579.745855734757,59217.71270161281,0.5556954971173874,0.5556954971173875),
    @("FAPs","Bgn","Fgfr3","FAPs",3,1,2511.398112333333,7534.194336999999,0.9109718135392577,0.9109718135392579,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,1431.116028649629,12880.04425784666,0.120865569949025,0.1208655699490251),
    @("FAPs","Bgn","Fgfr3","sCs",3,1,2511.398112333333,7534.194336999999,0.9109718135392577,0.9109718135392579,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,2775.55449998289,24979.99049984601,0.2344107464728452,0.2344107464728452),
    @("M2","Bgn","Fgfr3","ECs",3,1,5.092134333333333,15.276403,0.001847094980935658,0.001847094980935659,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,13.34115432039778,120.07038888358,0.001126733394380526,0.001126733394380526),
    @("M2","Bgn","Fgfr3","FAPs",3,1,5.092134333333333,15.276403,0.001847094980935658,0.001847094980935659,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,2.901744262959444,26.115698366635,0.0002450681616080002,0.0002450681616080003),
    @("M2","Bgn","Fgfr3","sCs",3,1,5.092134333333333,15.276403,0.001847094980935658,0.001847094980935659,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,5.627740298916334,50.649662690247,0.0004752934249471314,0.0004752934249471315),
    @("sCs","Bgn","Fgfr3","ECs",3,1,201.6251223333333,604.875367,0.07313647423921157,0.07313647423921157,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,528.2484113409579,4754.23570206862,0.04461346531883695,0.04461346531883695),
    @("sCs","Bgn","Fgfr3","FAPs",3,1,201.6251223333333,604.875367,0.07313647423921157,0.07313647423921157,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,114.8957399197794,1034.061659278015,0.009703573163961076,0.009703573163961076),
    @("sCs","Bgn","Fgfr3","sCs",3,1,201.6251223333333,604.875367,0.07313647423921157,0.07313647423921157,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,222.8326575757203,2005.493918181483,0.01881943575641354,0.01881943575641354)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
